$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update report header text (Volume/Number and date range) ---
# shared string run: "Volume " "30" "   Number  " "7"  ->  "...9"   (run 4 starts at char 21)
$ws.Range("C8").Characters(21, 1).Text = "9"
# shared string run: "Report Covering the Week  " "2/13/2023" "  Through  " "2/19/2023"
# run 2 ("2/13/2023") starts at char 27, len 9; run 4 ("2/19/2023") starts at char 47, len 9
$ws.Range("C9").Characters(27, 9).Text = "2/27/2023"
$ws.Range("C9").Characters(47, 9).Text = "3/5/2023"

# --- Crime statistics table updates ---
# Reference cell holding the canonical "text-right, General fmt" (style 14) look,
# used to restyle cells when a value reverts to a text placeholder ("0" / "***.*").
$refText = $ws.Range("C14")

# Row 15
$ws.Range("C15").Value = 1
$ws.Range("C15").NumberFormat = '#,##0'
$ws.Range("F15").Value = 1
$ws.Range("F15").NumberFormat = '#,##0'
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 2
$ws.Range("K15").Value = -50
$ws.Range("L15").Value = -33.333333333333
$ws.Range("M15").Value = 100
$ws.Range("N15").Value = -60
# Row 16
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = -75
$ws.Range("F16").Value = 7
$ws.Range("G16").Value = 15
$ws.Range("H16").Value = -53.333333333333
$ws.Range("I16").Value = 17
$ws.Range("J16").Value = 21
$ws.Range("K16").Value = -19.047619047619
$ws.Range("L16").Value = 54.545454545454
$ws.Range("M16").Value = -5.555555555555
$ws.Range("N16").Value = -83.962264150943
# Row 17
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 1
$ws.Range("D17").NumberFormat = '#,##0'
$ws.Range("E17").Value = 0
$ws.Range("E17").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F17").Value = 7
$ws.Range("G17").Value = 6
$ws.Range("H17").Value = 16.666666666666
$ws.Range("I17").Value = 15
$ws.Range("J17").Value = 12
$ws.Range("K17").Value = 25
$ws.Range("L17").Value = 25
$ws.Range("M17").Value = 7.142857142857
$ws.Range("N17").Value = 0
# Row 18
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 150
$ws.Range("F18").Value = 9
$ws.Range("G18").Value = 14
$ws.Range("H18").Value = -35.714285714285
$ws.Range("I18").Value = 32
$ws.Range("J18").Value = 26
$ws.Range("K18").Value = 23.076923076923
$ws.Range("L18").Value = 68.421052631578
$ws.Range("M18").Value = 3.225806451612
$ws.Range("N18").Value = -79.487179487179
# Row 19
$ws.Range("C19").Value = 17
$ws.Range("D19").Value = 16
$ws.Range("E19").Value = 6.25
$ws.Range("G19").Value = 62
$ws.Range("H19").Value = -3.225806451612
$ws.Range("I19").Value = 124
$ws.Range("J19").Value = 125
$ws.Range("K19").Value = -0.8
$ws.Range("L19").Value = 49.397590361445
$ws.Range("M19").Value = -8.148148148148
$ws.Range("N19").Value = -65.934065934065
# Row 20
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 2
$ws.Range("F20").Value = 11
$ws.Range("G20").Value = 5
$ws.Range("H20").Value = 120
$ws.Range("I20").Value = 17
$ws.Range("J20").Value = 8
$ws.Range("K20").Value = 112.5
$ws.Range("L20").Value = 183.333333333333
$ws.Range("M20").Value = 1600
$ws.Range("N20").Value = -91.981132075471
# Row 21
$ws.Range("C21").Value = 27
$ws.Range("D21").Value = 25
$ws.Range("E21").Value = 8
$ws.Range("F21").Value = 95
$ws.Range("G21").Value = 103
$ws.Range("H21").Value = -7.766990291262
$ws.Range("I21").Value = 208
$ws.Range("J21").Value = 196
$ws.Range("K21").Value = 6.122448979591
$ws.Range("L21").Value = 55.223880597014
$ws.Range("M21").Value = 4
$ws.Range("N21").Value = -75.757575757575
# Row 22
$ws.Range("C22").Value = "'0"
$refText.Copy() | Out-Null
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("D22").Value = 1
$ws.Range("D22").NumberFormat = '#,##0'
$ws.Range("E22").Value = -100
$ws.Range("E22").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("G22").Value = 4
$ws.Range("H22").Value = -75
$ws.Range("J22").Value = 6
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 100
$ws.Range("M22").Value = 0
# Row 23
$ws.Range("C23").Value = 3
$ws.Range("F23").Value = 6
$ws.Range("G23").Value = 1
$ws.Range("G23").NumberFormat = '#,##0'
$ws.Range("H23").Value = 500
$ws.Range("H23").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("I23").Value = 9
$ws.Range("J23").Value = 2
$ws.Range("K23").Value = 350
$ws.Range("L23").Value = 125
$ws.Range("M23").Value = 80
# Row 24
$ws.Range("C24").Value = 19
$ws.Range("D24").Value = 24
$ws.Range("E24").Value = -20.833333333333
$ws.Range("F24").Value = 73
$ws.Range("G24").Value = 103
$ws.Range("H24").Value = -29.126213592233
$ws.Range("I24").Value = 172
$ws.Range("J24").Value = 224
$ws.Range("K24").Value = -23.214285714285
$ws.Range("L24").Value = -38.351254480286
$ws.Range("M24").Value = 6.172839506172
# Row 25
$ws.Range("C25").Value = 4
$ws.Range("D25").Value = 3
$ws.Range("E25").Value = 33.333333333333
$ws.Range("F25").Value = 16
$ws.Range("G25").Value = 19
$ws.Range("H25").Value = -15.789473684210
$ws.Range("I25").Value = 36
$ws.Range("J25").Value = 35
$ws.Range("K25").Value = 2.857142857142
$ws.Range("L25").Value = 28.571428571428
$ws.Range("M25").Value = -2.702702702702
# Row 26
$ws.Range("C26").Value = 1
$ws.Range("C26").NumberFormat = '#,##0'
$ws.Range("F26").Value = 1
$ws.Range("F26").NumberFormat = '#,##0'
$ws.Range("G26").Value = 1
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 3
$ws.Range("K26").Value = -25
$ws.Range("L26").Value = -25
# Row 27
$ws.Range("C27").Value = 5
$ws.Range("C27").NumberFormat = '#,##0'
$ws.Range("D27").Value = 1
$ws.Range("D27").NumberFormat = '#,##0'
$ws.Range("E27").Value = 400
$ws.Range("E27").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F27").Value = 6
$ws.Range("F27").NumberFormat = '#,##0'
$ws.Range("H27").Value = 500
$ws.Range("I27").Value = 8
$ws.Range("J27").Value = 6
$ws.Range("K27").Value = 33.333333333333
$ws.Range("L27").Value = 166.666666666667
# Row 30
$ws.Range("D30").Value = "'0"
$refText.Copy() | Out-Null
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("E30").Value = "***.*"
